$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '68.159.86'
$c.ClearFormats()
$c = $ws.Range('E2')
$c.NumberFormat = '@'
$c.Value = '  +2.37%  '
$c.ClearFormats()

$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '2.535.44'
$c.ClearFormats()
$c = $ws.Range('E3')
$c.NumberFormat = '@'
$c.Value = '  +0.80%  '
$c.ClearFormats()

$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.ClearFormats()
$c = $ws.Range('E4')
$c.NumberFormat = '@'
$c.Value = '  -0.04%  '
$c.ClearFormats()

$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '594.46'
$c.ClearFormats()
$c = $ws.Range('E5')
$c.NumberFormat = '@'
$c.Value = '  +1.85%  '
$c.ClearFormats()

$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '177.30'
$c.ClearFormats()
$c = $ws.Range('E6')
$c.NumberFormat = '@'
$c.Value = '  +2.08%  '
$c.ClearFormats()

$c = $ws.Range('E7')
$c.NumberFormat = '@'
$c.Value = '  -0.11%  '
$c.ClearFormats()

$c = $ws.Range('E8')
$c.NumberFormat = '@'
$c.Value = '  +1.57%  '
$c.ClearFormats()

$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '2.532.94'
$c.ClearFormats()
$c = $ws.Range('E9')
$c.NumberFormat = '@'
$c.Value = '  +0.71%  '
$c.ClearFormats()

$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.142'
$c.ClearFormats()
$c = $ws.Range('E10')
$c.NumberFormat = '@'
$c.Value = '  +1.52%  '
$c.ClearFormats()

$c = $ws.Range('E11')
$c.NumberFormat = '@'
$c.Value = '  +3.04%  '
$c.ClearFormats()

$c = $ws.Range('E12')
$c.NumberFormat = '@'
$c.Value = '  +0.92%  '
$c.ClearFormats()

$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.348'
$c.ClearFormats()
$c = $ws.Range('E13')
$c.NumberFormat = '@'
$c.Value = '  -1.05%  '
$c.ClearFormats()

$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '27.00'
$c.ClearFormats()
$c = $ws.Range('E14')
$c.NumberFormat = '@'
$c.Value = '  +1.24%  '
$c.ClearFormats()

$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '2.988.61'
$c.ClearFormats()
$c = $ws.Range('E15')
$c.NumberFormat = '@'
$c.Value = '  +0.90%  '
$c.ClearFormats()

$c = $ws.Range('E16')
$c.NumberFormat = '@'
$c.Value = '  +1.38%  '
$c.ClearFormats()

$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '67.963.16'
$c.ClearFormats()
$c = $ws.Range('E17')
$c.NumberFormat = '@'
$c.Value = '  +2.40%  '
$c.ClearFormats()

$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '2.535.55'
$c.ClearFormats()
$c = $ws.Range('E18')
$c.NumberFormat = '@'
$c.Value = '  +1.00%  '
$c.ClearFormats()

$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '8.01'
$c.ClearFormats()
$c = $ws.Range('E19')
$c.NumberFormat = '@'
$c.Value = '  +3.55%  '
$c.ClearFormats()

$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '11.52'
$c.ClearFormats()
$c = $ws.Range('E20')
$c.NumberFormat = '@'
$c.Value = '  +2.06%  '
$c.ClearFormats()

$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '364.72'
$c.ClearFormats()
$c = $ws.Range('E21')
$c.NumberFormat = '@'
$c.Value = '  +4.67%  '
$c.ClearFormats()

$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '4.23'
$c.ClearFormats()
$c = $ws.Range('E22')
$c.NumberFormat = '@'
$c.Value = '  +0.37%  '
$c.ClearFormats()

$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '4.71'
$c.ClearFormats()
$c = $ws.Range('E23')
$c.NumberFormat = '@'
$c.Value = '  +1.75%  '
$c.ClearFormats()

$c = $ws.Range('E24')
$c.NumberFormat = '@'
$c.Value = '  -1.78%  '
$c.ClearFormats()

$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '71.04'
$c.ClearFormats()
$c = $ws.Range('E26')
$c.NumberFormat = '@'
$c.Value = '  +1.86%  '
$c.ClearFormats()

$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '10.27'
$c.ClearFormats()
$c = $ws.Range('E27')
$c.NumberFormat = '@'
$c.Value = '  +3.15%  '
$c.ClearFormats()

$c = $ws.Range('E28')
$c.NumberFormat = '@'
$c.Value = '  +0.60%  '
$c.ClearFormats()

$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '0.995'
$c.ClearFormats()

$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '0.0₃0993'
$c.ClearFormats()

$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '544.18'
$c.ClearFormats()
$c = $ws.Range('E31')
$c.NumberFormat = '@'
$c.Value = '  +2.92%  '
$c.ClearFormats()

$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '8.27'
$c.ClearFormats()
$c = $ws.Range('E32')
$c.NumberFormat = '@'
$c.Value = '  +1.76%  '
$c.ClearFormats()

$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '1.35'
$c.ClearFormats()
$c = $ws.Range('E33')
$c.NumberFormat = '@'
$c.Value = '  +1.57%  '
$c.ClearFormats()

$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '1.88'
$c.ClearFormats()
$c = $ws.Range('E34')
$c.NumberFormat = '@'
$c.Value = '  +2.17%  '
$c.ClearFormats()

$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.ClearFormats()
$c = $ws.Range('E36')
$c.NumberFormat = '@'
$c.Value = '  -0.02%  '
$c.ClearFormats()

$c = $ws.Range('E37')
$c.NumberFormat = '@'
$c.Value = '  +0.91%  '
$c.ClearFormats()

$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '156.77'
$c.ClearFormats()
$c = $ws.Range('E38')
$c.NumberFormat = '@'
$c.Value = '  +0.71%  '
$c.ClearFormats()

$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '18.87'
$c.ClearFormats()
$c = $ws.Range('E39')
$c.NumberFormat = '@'
$c.Value = '  +1.37%  '
$c.ClearFormats()

$c = $ws.Range('E40')
$c.NumberFormat = '@'
$c.Value = '  +1.62%  '
$c.ClearFormats()

$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.358'
$c.ClearFormats()
$c = $ws.Range('E41')
$c.NumberFormat = '@'
$c.Value = '  +0.79%  '
$c.ClearFormats()

$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '5.22'
$c.ClearFormats()
$c = $ws.Range('E42')
$c.NumberFormat = '@'
$c.Value = '  +2.25%  '
$c.ClearFormats()

$c = $ws.Range('E43')
$c.NumberFormat = '@'
$c.Value = '  +0.97%  '
$c.ClearFormats()

$c = $ws.Range('B44')
$c.NumberFormat = '@'
$c.Value = 'dogwifhat'
$c.ClearFormats()
$c = $ws.Range('C44')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$c.ClearFormats()
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '2.53'
$c.ClearFormats()
$c = $ws.Range('E44')
$c.NumberFormat = '@'
$c.Value = '  -0.39%  '
$c.ClearFormats()

$c = $ws.Range('B45')
$c.NumberFormat = '@'
$c.Value = 'USDe'
$c.ClearFormats()
$c = $ws.Range('C45')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$c.ClearFormats()
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.ClearFormats()
$c = $ws.Range('E45')
$c.NumberFormat = '@'
$c.Value = '  +0.00%  '
$c.ClearFormats()

$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.563'
$c.ClearFormats()
$c = $ws.Range('E46')
$c.NumberFormat = '@'
$c.Value = '  +0.80%  '
$c.ClearFormats()

$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '147.42'
$c.ClearFormats()
$c = $ws.Range('E47')
$c.NumberFormat = '@'
$c.Value = '  -0.44%  '
$c.ClearFormats()

$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '3.74'
$c.ClearFormats()
$c = $ws.Range('E48')
$c.NumberFormat = '@'
$c.Value = '  +1.66%  '
$c.ClearFormats()

$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.0₆0279'
$c.ClearFormats()
$c = $ws.Range('E49')
$c.NumberFormat = '@'
$c.Value = '  +1.89%  '
$c.ClearFormats()

$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '1.71'
$c.ClearFormats()
$c = $ws.Range('E50')
$c.NumberFormat = '@'
$c.Value = '  -1.50%  '
$c.ClearFormats()

$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.0759'
$c.ClearFormats()
$c = $ws.Range('E51')
$c.NumberFormat = '@'
$c.Value = '  +1.39%  '
$c.ClearFormats()
